# Commit: "add the NA's under duplicate_image_filename"
#
# Column E on the sheet is headed "duplicate_image_filename" (E1).
# For every word-trial row (the practice rows 2-5 and the generic /
# unique_video / unique_audio word rows 6-21) that column was blank;
# this fills those cells in with the literal text "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
